# Trade #122 (MarketMaking) closes with an early exit, and a new Trade #151
# (momentum, OPEN) is recorded. This updates the roll-up sheets (Summary,
# Strategy Status), the master "All Trades" log, the strategy-specific
# "momentum" log (new trade) and the strategy-specific "MarketMaking" log
# (closed trade).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1499.45              # Current Capital
$wsSummary.Range("B4").Value = 0.5600000000000001   # Total P&L $
$wsSummary.Range("B6").Value = 122                  # Total Trades
$wsSummary.Range("B8").Value = 43                   # Losing Trades
$wsSummary.Range("B9").Value = 47.54                # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 6)
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C6").Value = 99.53    # Capital
$wsStatus.Range("D6").Value = 42       # Trades
$wsStatus.Range("E6").Value = -0.28    # P&L $
$wsStatus.Range("F6").Value = -0.47    # P&L %
$wsStatus.Range("G6").Value = 47.62    # Win Rate %

# ---------------------------------------------------------------------
# All Trades sheet - close out Trade #122 (row 123)
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")
$wsAll.Cells.Item(123, 7).Value = 0.001             # G123 Exit Price
$wsAll.Cells.Item(123, 8).Value = "CLOSED"          # H123 Status
$wsAll.Cells.Item(123, 9).Value = -90               # I123 P&L %
$wsAll.Cells.Item(123, 10).Value = -0.01            # J123 P&L $
$wsAll.Cells.Item(123, 11).Value = 99.53            # K123 Capital After
$wsAll.Cells.Item(123, 12).Value = "early_exit"     # L123 Exit Reason
$wsAll.Cells.Item(123, 13).Value = 0.16             # M123 Duration (min)

# All Trades sheet - append new Trade #151 (row 152, momentum, still OPEN)
$wsAll.Cells.Item(152, 1).Value = 151               # A152 Trade #

$wsAll.Cells.Item(152, 2).NumberFormat = "@"        # force text so the
$wsAll.Cells.Item(152, 2).Value = "2026-02-18"      # date-like string isn't
$wsAll.Cells.Item(152, 2).ClearFormats()            # auto-converted to a date

$wsAll.Cells.Item(152, 3).Value = "00:29:42"        # C152 Time
$wsAll.Cells.Item(152, 4).Value = "momentum"        # D152 Strategy
$wsAll.Cells.Item(152, 5).Value = "UP"              # E152 Side
$wsAll.Cells.Item(152, 6).Value = 0.01               # F152 Entry Price
$wsAll.Cells.Item(152, 8).Value = "OPEN"            # H152 Status
$wsAll.Cells.Item(152, 9).Value = 0                 # I152 P&L %
$wsAll.Cells.Item(152, 10).Value = 0                # J152 P&L $
$wsAll.Cells.Item(152, 11).Value = 99.23374292899115 # K152 Capital After
$wsAll.Cells.Item(152, 13).Value = 0                # M152 Duration (min)
$wsAll.Cells.Item(152, 14).Value = 0                # N152 Entry Slippage (bps)
$wsAll.Cells.Item(152, 15).Value = 0                # O152 Exit Slippage (bps)
$wsAll.Cells.Item(152, 16).Value = 0.9              # P152 Confidence
$wsAll.Cells.Item(152, 17).Value = "Upward momentum: 2.941% over 10 samples" # Q152 Entry Reason

# ---------------------------------------------------------------------
# momentum sheet - append new Trade #151 (row 37, still OPEN)
# ---------------------------------------------------------------------
$wsMomentum = $wb.Worksheets.Item("momentum")
$wsMomentum.Cells.Item(37, 1).Value = 151           # A37 Trade #

$wsMomentum.Cells.Item(37, 2).NumberFormat = "@"
$wsMomentum.Cells.Item(37, 2).Value = "2026-02-18"  # B37 Date
$wsMomentum.Cells.Item(37, 2).ClearFormats()

$wsMomentum.Cells.Item(37, 3).Value = "00:29:42"    # C37 Time
$wsMomentum.Cells.Item(37, 4).Value = "momentum"    # D37 Strategy
$wsMomentum.Cells.Item(37, 5).Value = "UP"          # E37 Side
$wsMomentum.Cells.Item(37, 6).Value = 0.01           # F37 Entry Price
$wsMomentum.Cells.Item(37, 8).Value = "OPEN"        # H37 Status
$wsMomentum.Cells.Item(37, 9).Value = 0             # I37 P&L %
$wsMomentum.Cells.Item(37, 10).Value = 0            # J37 P&L $
$wsMomentum.Cells.Item(37, 11).Value = 99.23374292899115 # K37 Capital After
$wsMomentum.Cells.Item(37, 12).Value = 0            # L37 Entry Slippage (bps)
$wsMomentum.Cells.Item(37, 13).Value = 0            # M37 Exit Slippage (bps)
$wsMomentum.Cells.Item(37, 14).Value = 0.9          # N37 Confidence
$wsMomentum.Cells.Item(37, 15).Value = "Upward momentum: 2.941% over 10 samples" # O37 Entry Reason
$wsMomentum.Cells.Item(37, 17).Value = 0            # Q37 Duration (min)

# ---------------------------------------------------------------------
# MarketMaking sheet - close out Trade #122 (row 43)
# ---------------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")
$wsMM.Cells.Item(43, 7).Value = 0.001               # G43 Exit Price
$wsMM.Cells.Item(43, 8).Value = "CLOSED"            # H43 Status
$wsMM.Cells.Item(43, 9).Value = -90                 # I43 P&L %
$wsMM.Cells.Item(43, 10).Value = -0.01              # J43 P&L $
$wsMM.Cells.Item(43, 11).Value = 99.53              # K43 Capital After
$wsMM.Cells.Item(43, 16).Value = "early_exit"       # P43 Exit Reason
$wsMM.Cells.Item(43, 17).Value = 0.16               # Q43 Duration (min)
